# Apply cryptocurrency price/volume/name updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($rng, $val)
    # Force the cell to stay a text value (avoid Excel coercing numeric-looking
    # strings like "309.57" or "1.00" into real numbers), then restore the
    # original "Normal" cell style so no stray number-format/style is left behind.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '43.332.14'
Set-TextValue $ws.Range("E2") '  +1.69%  '
Set-TextValue $ws.Range("D3") '2.347.13'
Set-TextValue $ws.Range("E3") '  +3.29%  '
Set-TextValue $ws.Range("E4") '  +0.02%  '
Set-TextValue $ws.Range("D5") '309.57'
Set-TextValue $ws.Range("E5") '  +0.40%  '
Set-TextValue $ws.Range("D6") '104.22'
Set-TextValue $ws.Range("E6") '  +4.80%  '
Set-TextValue $ws.Range("D7") '0.525'
Set-TextValue $ws.Range("E7") '  -0.29%  '
Set-TextValue $ws.Range("E8") '  -0.08%  '
Set-TextValue $ws.Range("D9") '0.521'
Set-TextValue $ws.Range("E9") '  +4.71%  '
Set-TextValue $ws.Range("D10") '36.12'
Set-TextValue $ws.Range("E10") '  +0.93%  '
Set-TextValue $ws.Range("D11") '52.70'
Set-TextValue $ws.Range("E11") '  +1.50%  '
Set-TextValue $ws.Range("D12") '0.0812'
Set-TextValue $ws.Range("E12") '  -0.59%  '
Set-TextValue $ws.Range("E13") '  -1.49%  '
Set-TextValue $ws.Range("D14") '6.98'
Set-TextValue $ws.Range("E14") '  +2.88%  '
Set-TextValue $ws.Range("D15") '2.708.44'
Set-TextValue $ws.Range("E15") '  +3.21%  '
Set-TextValue $ws.Range("D16") '15.63'
Set-TextValue $ws.Range("E16") '  +6.72%  '
Set-TextValue $ws.Range("D17") '2.343.75'
Set-TextValue $ws.Range("E17") '  +0.44%  '
Set-TextValue $ws.Range("D18") '0.807'
Set-TextValue $ws.Range("E18") '  +2.13%  '
Set-TextValue $ws.Range("D19") '43.302.74'
Set-TextValue $ws.Range("E19") '  +1.84%  '
$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D20") '11.99'
Set-TextValue $ws.Range("E20") '  -2.68%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D21") '0.0₃0925'
Set-TextValue $ws.Range("E21") '  +2.01%  '
Set-TextValue $ws.Range("D22") '6.28'
Set-TextValue $ws.Range("E22") '  +4.79%  '
Set-TextValue $ws.Range("D23") '68.18'
Set-TextValue $ws.Range("E23") '  +1.06%  '
Set-TextValue $ws.Range("D24") '241.75'
Set-TextValue $ws.Range("E24") '  +1.87%  '
Set-TextValue $ws.Range("E25") '  +3.54%  '
Set-TextValue $ws.Range("E26") '  +0.76%  '
Set-TextValue $ws.Range("D27") '1.00'
Set-TextValue $ws.Range("E27") '  -0.48%  '
Set-TextValue $ws.Range("D28") '25.36'
Set-TextValue $ws.Range("E28") '  +7.40%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D29") '2.22'
Set-TextValue $ws.Range("E29") '  +3.46%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D30") '36.53'
Set-TextValue $ws.Range("E30") '  -5.00%  '
Set-TextValue $ws.Range("D31") '9.59'
Set-TextValue $ws.Range("E31") '  +0.13%  '
Set-TextValue $ws.Range("D32") '162.15'
Set-TextValue $ws.Range("E33") '  +0.61%  '
Set-TextValue $ws.Range("E34") '  +0.01%  '
Set-TextValue $ws.Range("D35") '18.28'
Set-TextValue $ws.Range("E35") '  +3.13%  '
Set-TextValue $ws.Range("E36") '  +6.69%  '
Set-TextValue $ws.Range("D37") '3.10'
Set-TextValue $ws.Range("E37") '  +0.56%  '
Set-TextValue $ws.Range("E38") '  +1.39%  '
Set-TextValue $ws.Range("D39") '4.58'
Set-TextValue $ws.Range("E39") '  +10.34%  '
Set-TextValue $ws.Range("D40") '1.90'
Set-TextValue $ws.Range("E40") '  +5.44%  '
Set-TextValue $ws.Range("D41") '0.106'
Set-TextValue $ws.Range("E41") '  +2.63%  '
Set-TextValue $ws.Range("E42") '  +0.15%  '
Set-TextValue $ws.Range("D43") '2.37'
Set-TextValue $ws.Range("E43") '  +4.59%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D44") '19.99'
Set-TextValue $ws.Range("E44") '  +4.93%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D45") '0.0291'
Set-TextValue $ws.Range("E45") '  +2.05%  '
Set-TextValue $ws.Range("D46") '1.981.96'
Set-TextValue $ws.Range("E46") '  +1.83%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D47") '10.42'
Set-TextValue $ws.Range("E47") '  +7.00%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D48") '3.06'
Set-TextValue $ws.Range("E48") '  +4.19%  '
Set-TextValue $ws.Range("D49") '58.87'
Set-TextValue $ws.Range("E49") '  +8.35%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D50") '1.59'
Set-TextValue $ws.Range("E50") '  +6.29%  '
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range("D51") '2.90'
Set-TextValue $ws.Range("E51") '  -0.70%  '
